$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being updated, in order, for each row (2..25)
$cols = @("B", "C", "D", "E", "F", "I")

# New values per row (row 2 first, row 25 last), matching $cols order
$data = @(
    @(15.85394732810968,13.48514581944743,4.969206946898082,16.49250889861493,31.91405166130836,23.19874287244302),
    @(15.15797451134646,12.68701410148378,4.982433563494094,15.54982305427502,31.27412109795081,23.02195282445205),
    @(14.71890366318364,12.17409127707578,4.991782024880804,14.94723132536202,30.88727418257572,22.92057014189517),
    @(14.53732229752811,11.95951077681816,4.995895811352505,14.69596131537296,30.73137712178652,22.88107784963083),
    @(14.50701978616817,11.92355038218285,4.996597147649619,14.6539022716624,30.70560235488173,22.87463072286543),
    @(14.7164651308351,12.17121960191544,4.991836278894994,14.94386532916705,30.88516434215108,22.92003013485543),
    @(15.61657884908547,13.21479768401042,4.973510122851643,16.17255248081426,31.69227508032762,23.1363112085674),
    @(17.27736391435519,15.07334164428957,4.947507018560797,18.4744944591479,33.31355011951186,23.61624884215468),
    @(18.42130341069321,16.31769835030382,4.934728174615995,20.14109724497979,34.51542795369839,24.00121645478992),
    @(18.92308297424428,16.85658083344993,4.930345928834477,20.85794621852676,35.06192119520299,24.18294992050313),
    @(19.11027981864848,17.06307154061995,4.928896959247793,21.12350304161436,35.26862224612128,24.25267804753229),
    @(19.07009088930396,17.01734708079981,4.929199584005914,21.06657231982961,35.22411954578326,24.23762106804894),
    @(18.93854094049538,16.87312282559413,4.930222475658698,20.87991171968907,35.07893266580561,24.1886685079885),
    @(18.85759224448678,16.78645957656379,4.930876578551215,20.764809766483,34.98996370356716,24.15880085614728),
    @(18.38812401053178,16.28192909908904,4.935043762929912,20.09342190971258,34.47969031268256,23.98946931219963),
    @(18.0952479518424,15.96541014474897,4.937970056798055,19.6709986320909,34.16643939399061,23.8872545153566),
    @(17.92504916041692,15.78080061781358,4.939787586385396,19.42414247111821,33.9862544000481,23.82908802385981),
    @(17.86712786630715,15.71785789745826,4.940425910077389,19.33989197422478,33.92525142162285,23.80950235699019),
    @(18.12660672119382,15.99936912980512,4.93764460334295,19.71636856486098,34.1997882059074,23.89807108356061),
    @(18.97725779326987,16.91453994987527,4.929916277313257,20.9348982220559,35.12158580854742,24.20302271226781),
    @(19.51674223916369,17.5234870285506,4.926094454991147,21.69691642342698,35.72252169950087,24.40760760245359),
    @(19.23035676717818,17.19946751425355,4.928020215389134,21.29334379216612,35.40199552985129,24.29794763539691),
    @(18.1124350809237,15.98402448686858,4.937791320202155,19.69586927850612,34.18471148136781,23.89317904687864),
    @(16.84067190314043,14.59145038599028,4.95345016798223,17.8238655196067,32.87216146039342,23.48059623871805)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $values[$j]
    }
}
